$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 85, shifting rows 85:120 down to 86:121
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with its new data
$ws.Range("A85").Value = 7
$ws.Range("B85").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C85").Value = "Ñuble"
$ws.Range("D85").Value = 44992
$ws.Range("E85").Value = 16
$ws.Range("F85").Value = 100112030
$ws.Range("G85").Value = "Poroto granado"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 50
$ws.Range("K85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = 30000
$ws.Range("N85").Value = "$/saco 25 kilos"
$ws.Range("O85").Value = "Provincia de Diguillín"
$ws.Range("P85").Value = 1200
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"
